# Apply the update: swap several adjacent match rows back to their correct
# chronological order (columns F:V only - "home" through "url_partida"),
# and append 4 new match rows (165-168) for matches played on
# 30/31 Oct - 06/07 Nov 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param([int]$Row1, [int]$Row2)

    $rng1 = $ws.Range("F$Row1" + ":V$Row1")
    $rng2 = $ws.Range("F$Row2" + ":V$Row2")

    $v1 = $rng1.Value2
    $v2 = $rng2.Value2

    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Pairs of rows whose F:V (home..url_partida) contents were swapped
Swap-MatchRows 60 61
Swap-MatchRows 64 65
Swap-MatchRows 80 81
Swap-MatchRows 82 83
Swap-MatchRows 98 99
Swap-MatchRows 102 103
Swap-MatchRows 110 111
Swap-MatchRows 125 126
Swap-MatchRows 140 141

function Add-MatchRow {
    param([int]$Row, [object[]]$Values)

    # Copy formatting (styles) from the previous row so the new row gets the
    # same bold/centered index style (col A) and date-time style (col E).
    $prevRow = $Row - 1
    $ws.Range("A$prevRow" + ":V$prevRow").Copy()
    $ws.Range("A$Row" + ":V$Row").PasteSpecial(-4122) | Out-Null

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value2 = $Values[$i]
    }

    # Column D ("temporada") must stay a text value ("2023") rather than be
    # auto-coerced to a number. Force text formatting, set the value, then
    # restore the plain/default cell style (copied from a known text cell)
    # so no stray numeric style is introduced.
    $dCell = $ws.Cells.Item($Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value2 = "2023"
    $ws.Range("D160").Copy()
    $dCell.PasteSpecial(-4122) | Out-Null
}

Add-MatchRow 165 @(164, "argentina", "copa-de-la-liga-profesional", "2023", 45236.9375, "Barracas Central", 1, "Rosario Central", 1, 2.47, "31/10/2023 01:12", 2.37, "06/11/2023 22:25", 3.12, "31/10/2023 01:12", 2.99, "06/11/2023 22:25", 3.15, "31/10/2023 01:12", 3.61, "06/11/2023 22:25", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-rosario-central/rmfjKtj2/")
Add-MatchRow 166 @(165, "argentina", "copa-de-la-liga-profesional", "2023", 45236.9375, "Godoy Cruz", 2, "Platense", 0, 1.98, "31/10/2023 23:12", 1.87, "06/11/2023 22:02", 3.2, "31/10/2023 23:12", 3.2, "06/11/2023 22:02", 4.43, "31/10/2023 23:12", 5.29, "06/11/2023 22:02", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-platense/IHAQLgR6/")
Add-MatchRow 167 @(166, "argentina", "copa-de-la-liga-profesional", "2023", 45237.04166666666, "Banfield", 2, "Colon Santa Fe", 1, 2.01, "30/10/2023 20:12", 2.23, "07/11/2023 00:57", 3.25, "30/10/2023 20:12", 3, "07/11/2023 00:57", 4.19, "30/10/2023 20:12", 3.98, "07/11/2023 00:46", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/banfield-colon-santa-fe/CbKBROcS/")
Add-MatchRow 168 @(167, "argentina", "copa-de-la-liga-profesional", "2023", 45237.04166666666, "Belgrano", 0, "Tigre", 3, 2.22, "31/10/2023 23:12", 2.57, "07/11/2023 00:59", 3.05, "31/10/2023 23:12", 2.79, "07/11/2023 00:59", 3.53, "31/10/2023 23:12", 3.5, "07/11/2023 00:59", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/ca-belgrano-de-cordoba-tigre/bexwrpr9/")
